$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '35.114.10'
$ws.Cells.Item(2, 5).Value = '  +0.70%  '

$ws.Cells.Item(3, 4).Value = '1.855.91'
$ws.Cells.Item(3, 5).Value = '  +1.60%  '

$ws.Cells.Item(4, 5).Value = '  +0.40%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '238.73'
$ws.Cells.Item(5, 5).Value = '  +3.57%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.622'
$ws.Cells.Item(6, 5).Value = '  +0.97%  '

$ws.Cells.Item(7, 5).Value = '  +0.38%  '

$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '41.97'
$ws.Cells.Item(8, 5).Value = '  +5.08%  '

$ws.Cells.Item(9, 5).Value = '  +2.54%  '

$ws.Cells.Item(10, 5).Value = '  +1.62%  '

$ws.Cells.Item(11, 5).Value = '  +0.01%  '

$ws.Cells.Item(12, 4).Value = '2.123.23'
$ws.Cells.Item(12, 5).Value = '  +1.46%  '

$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '11.47'
$ws.Cells.Item(13, 5).Value = '  +1.68%  '

$ws.Cells.Item(14, 4).Value = '1.863.91'
$ws.Cells.Item(14, 5).Value = '  +2.02%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.676'
$ws.Cells.Item(15, 5).Value = '  +1.38%  '

$ws.Cells.Item(16, 5).Value = '  +1.59%  '

$ws.Cells.Item(17, 4).Value = '35.092.12'
$ws.Cells.Item(17, 5).Value = '  +0.57%  '

$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '69.90'
$ws.Cells.Item(18, 5).Value = '  +0.61%  '

$ws.Cells.Item(19, 4).Value = '0.0₃0794'
$ws.Cells.Item(19, 5).Value = '  +1.27%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '240.75'
$ws.Cells.Item(20, 5).Value = '  +0.53%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '12.25'
$ws.Cells.Item(21, 5).Value = '  +1.02%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '4.71'
$ws.Cells.Item(22, 5).Value = '  +1.42%  '

$ws.Cells.Item(23, 5).Value = '  +0.31%  '

$ws.Cells.Item(24, 5).Value = '  +0.43%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '167.70'
$ws.Cells.Item(25, 5).Value = '  -3.34%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '1.88'
$ws.Cells.Item(26, 5).Value = '  +24.93%  '

$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '7.96'
$ws.Cells.Item(27, 5).Value = '  +3.27%  '

$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '17.65'
$ws.Cells.Item(28, 5).Value = '  +1.84%  '

$ws.Cells.Item(29, 5).Value = '  +0.06%  '

$ws.Cells.Item(30, 5).Value = '  +0.38%  '

$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '0.0557'
$ws.Cells.Item(31, 5).Value = '  +1.36%  '

$ws.Cells.Item(32, 5).Value = '  +2.07%  '

$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '1.81'
$ws.Cells.Item(33, 5).Value = '  +27.35%  '

$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '4.01'
$ws.Cells.Item(34, 5).Value = '  +2.34%  '

$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.837'
$ws.Cells.Item(35, 5).Value = '  +19.89%  '

$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '2.03'
$ws.Cells.Item(36, 5).Value = '  +11.75%  '

$ws.Cells.Item(37, 5).Value = '  +6.90%  '

$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '1.09'
$ws.Cells.Item(38, 5).Value = '  +6.93%  '

$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '90.55'
$ws.Cells.Item(39, 5).Value = '  -1.59%  '

$ws.Cells.Item(40, 5).Value = '  +4.06%  '

$ws.Cells.Item(41, 4).Value = '1.340.86'
$ws.Cells.Item(41, 5).Value = '  +0.20%  '

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '14.92'
$ws.Cells.Item(42, 5).Value = '  +3.15%  '

$ws.Cells.Item(43, 5).Value = '  +3.65%  '

$ws.Cells.Item(44, 2).Value = 'Gas'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '12.55'
$ws.Cells.Item(44, 5).Value = '  +45.77%  '

$ws.Cells.Item(45, 2).Value = 'HuobiToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '2.41'
$ws.Cells.Item(45, 5).Value = '  -0.46%  '

$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.0557'
$ws.Cells.Item(46, 5).Value = '  +6.74%  '

$ws.Cells.Item(47, 5).Value = '  -0.14%  '

$ws.Cells.Item(48, 5).Value = '  +5.36%  '

$ws.Cells.Item(49, 4).Value = '2.038.88'
$ws.Cells.Item(49, 5).Value = '  +1.44%  '

$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.0680'
$ws.Cells.Item(50, 5).Value = '  +1.74%  '

$ws.Cells.Item(51, 5).Value = '  +0.43%  '
